$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 171950
$ws.Range("C4").Value = 162746
$ws.Range("C5").Value = 9204
$ws.Range("C8").Value = 66
